# fixed typo: list_name -> 'list name'
# (XlsForm struct changed to accomodate future features)

$wb = $excel.ActiveWorkbook
$choices = $wb.Worksheets.Item("choices")

# Fix the typo in the "choices" sheet header cell (A1):
# "list_name" -> "list name"
$choices.Range("A1").Value = "list name"

# The "choices" sheet is now the active/selected tab (it was "survey" before).
[void]$choices.Select()
[void]$choices.Range("A1").Select()

$wb.Save()
